# 1 camada 15 neu
# Update metrics rows 2-26 with new model_9_7_* values (models reindexed/reordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_9_7_0"
$ws.Range("B2").Value = 0.7587179157229249
$ws.Range("C2").Value = 0.7799532191184506
$ws.Range("D2").Value = 0.9732464376637572
$ws.Range("E2").Value = 0.9141328234428118
$ws.Range("F2").Value = 0.2670281827449799
$ws.Range("G2").Value = 0.2104353606700897
$ws.Range("H2").Value = 0.06531330198049545
$ws.Range("I2").Value = 0.1421425491571426
$ws.Range("A3").Value = "model_9_7_1"
$ws.Range("B3").Value = 0.7972219277606299
$ws.Range("C3").Value = 0.7492066339061014
$ws.Range("D3").Value = 0.9605046115323883
$ws.Range("E3").Value = 0.8958862268751013
$ws.Range("F3").Value = 0.2244155704975128
$ws.Range("G3").Value = 0.2398389726877213
$ws.Range("H3").Value = 0.09641984105110168
$ws.Range("I3").Value = 0.1723475605249405
$ws.Range("A4").Value = "model_9_7_2"
$ws.Range("B4").Value = 0.80190298743632
$ws.Range("C4").Value = 0.7416314807845422
$ws.Range("D4").Value = 0.9611490682799176
$ws.Range("E4").Value = 0.8940165679359064
$ws.Range("F4").Value = 0.2192350178956985
$ws.Range("G4").Value = 0.2470832467079163
$ws.Range("H4").Value = 0.09484653174877167
$ws.Range("I4").Value = 0.1754425764083862
$ws.Range("A5").Value = "model_9_7_3"
$ws.Range("B5").Value = 0.8601416534317192
$ws.Range("C5").Value = 0.7733421478147926
$ws.Range("D5").Value = 0.9538995397028258
$ws.Range("E5").Value = 0.8986839488134108
$ws.Range("F5").Value = 0.1547819674015045
$ws.Range("G5").Value = 0.216757670044899
$ws.Range("H5").Value = 0.1125447601079941
$ws.Range("I5").Value = 0.1677162796258926
$ws.Range("A6").Value = "model_9_7_4"
$ws.Range("B6").Value = 0.8794811137173297
$ws.Range("C6").Value = 0.732027535627271
$ws.Range("D6").Value = 0.9184805961116117
$ws.Range("E6").Value = 0.8614670457287493
$ws.Range("F6").Value = 0.1333788931369781
$ws.Range("G6").Value = 0.2562676668167114
$ws.Range("H6").Value = 0.199012815952301
$ws.Range("I6").Value = 0.2293242961168289
$ws.Range("A7").Value = "model_9_7_5"
$ws.Range("B7").Value = 0.8910793689412591
$ws.Range("C7").Value = 0.7261107164228393
$ws.Range("D7").Value = 0.9246496271479849
$ws.Range("E7").Value = 0.8639387384512229
$ws.Range("F7").Value = 0.120543047785759
$ws.Range("G7").Value = 0.2619260549545288
$ws.Range("H7").Value = 0.1839523911476135
$ws.Range("I7").Value = 0.225232720375061
$ws.Range("A8").Value = "model_9_7_8"
$ws.Range("B8").Value = 0.897164092096124
$ws.Range("C8").Value = 0.6979406884373178
$ws.Range("D8").Value = 0.9338625567694548
$ws.Range("E8").Value = 0.86171703777787
$ws.Range("F8").Value = 0.1138090416789055
$ws.Range("G8").Value = 0.2888656854629517
$ws.Range("H8").Value = 0.1614609211683273
$ws.Range("I8").Value = 0.2289104610681534
$ws.Range("A9").Value = "model_9_7_6"
$ws.Range("B9").Value = 0.8978966253786839
$ws.Range("C9").Value = 0.7259137755971721
$ws.Range("D9").Value = 0.930289623309864
$ws.Range("E9").Value = 0.8677926756423322
$ws.Range("F9").Value = 0.1129983514547348
$ws.Range("G9").Value = 0.262114405632019
$ws.Range("H9").Value = 0.1701835095882416
$ws.Range("I9").Value = 0.2188529968261719
$ws.Range("A10").Value = "model_9_7_9"
$ws.Range("B10").Value = 0.8983847029090337
$ws.Range("C10").Value = 0.7039491684752335
$ws.Range("D10").Value = 0.9334645774793513
$ws.Range("E10").Value = 0.8632780183246735
$ws.Range("F10").Value = 0.1124581918120384
$ws.Range("G10").Value = 0.2831196188926697
$ws.Range("H10").Value = 0.1624325215816498
$ws.Range("I10").Value = 0.2263264507055283
$ws.Range("A11").Value = "model_9_7_7"
$ws.Range("B11").Value = 0.8989599486019595
$ws.Range("C11").Value = 0.7281484259033211
$ws.Range("D11").Value = 0.9301784621590716
$ws.Range("E11").Value = 0.8683991267864132
$ws.Range("F11").Value = 0.1118215695023537
$ws.Range("G11").Value = 0.2599773705005646
$ws.Range("H11").Value = 0.1704548746347427
$ws.Range("I11").Value = 0.21784907579422
$ws.Range("A12").Value = "model_9_7_10"
$ws.Range("B12").Value = 0.8995707380421577
$ws.Range("C12").Value = 0.709305896173531
$ws.Range("D12").Value = 0.9339001531009218
$ws.Range("E12").Value = 0.8652195840683823
$ws.Range("F12").Value = 0.1111456006765366
$ws.Range("G12").Value = 0.2779968976974487
$ws.Range("H12").Value = 0.1613691449165344
$ws.Range("I12").Value = 0.2231124192476273
$ws.Range("A13").Value = "model_9_7_11"
$ws.Range("B13").Value = 0.9001124652884568
$ws.Range("C13").Value = 0.7109293644357265
$ws.Range("D13").Value = 0.9341121345928806
$ws.Range("E13").Value = 0.8658626917567929
$ws.Range("F13").Value = 0.1105460673570633
$ws.Range("G13").Value = 0.2764443159103394
$ws.Range("H13").Value = 0.1608516275882721
$ws.Range("I13").Value = 0.2220478504896164
$ws.Range("A14").Value = "model_9_7_15"
$ws.Range("B14").Value = 0.9003235094563341
$ws.Range("C14").Value = 0.7090279062858289
$ws.Range("D14").Value = 0.9334002673737887
$ws.Range("E14").Value = 0.8647872894507644
$ws.Range("F14").Value = 0.1103125065565109
$ws.Range("G14").Value = 0.278262734413147
$ws.Range("H14").Value = 0.1625895202159882
$ws.Range("I14").Value = 0.2238280475139618
$ws.Range("A15").Value = "model_9_7_12"
$ws.Range("B15").Value = 0.9003651108087392
$ws.Range("C15").Value = 0.7108636695988717
$ws.Range("D15").Value = 0.9334874996331853
$ws.Range("E15").Value = 0.8654088955074073
$ws.Range("F15").Value = 0.1102664619684219
$ws.Range("G15").Value = 0.276507169008255
$ws.Range("H15").Value = 0.1623765528202057
$ws.Range("I15").Value = 0.2227990329265594
$ws.Range("A16").Value = "model_9_7_14"
$ws.Range("B16").Value = 0.9005646039232871
$ws.Range("C16").Value = 0.7114466442808289
$ws.Range("D16").Value = 0.9332612501194234
$ws.Range("E16").Value = 0.8654301832605388
$ws.Range("F16").Value = 0.1100456863641739
$ws.Range("G16").Value = 0.275949627161026
$ws.Range("H16").Value = 0.162928894162178
$ws.Range("I16").Value = 0.222763791680336
$ws.Range("A17").Value = "model_9_7_13"
$ws.Range("B17").Value = 0.9006682963766114
$ws.Range("C17").Value = 0.7125192790214938
$ws.Range("D17").Value = 0.9333748707263577
$ws.Range("E17").Value = 0.8658369498351863
$ws.Range("F17").Value = 0.1099309250712395
$ws.Range("G17").Value = 0.2749238610267639
$ws.Range("H17").Value = 0.1626515090465546
$ws.Range("I17").Value = 0.2220904231071472
$ws.Range("A18").Value = "model_9_7_16"
$ws.Range("B18").Value = 0.9009976993102151
$ws.Range("C18").Value = 0.712578044845916
$ws.Range("D18").Value = 0.9339619399510605
$ws.Range("E18").Value = 0.8662633656064305
$ws.Range("F18").Value = 0.1095663756132126
$ws.Range("G18").Value = 0.2748676538467407
$ws.Range("H18").Value = 0.1612183153629303
$ws.Range("I18").Value = 0.2213845700025558
$ws.Range("A19").Value = "model_9_7_17"
$ws.Range("B19").Value = 0.9019957757101483
$ws.Range("C19").Value = 0.7205050246400748
$ws.Range("D19").Value = 0.9332487798237872
$ws.Range("E19").Value = 0.8681917009823001
$ws.Range("F19").Value = 0.1084617972373962
$ws.Range("G19").Value = 0.2672868967056274
$ws.Range("H19").Value = 0.1629593372344971
$ws.Range("I19").Value = 0.218192458152771
$ws.Range("A20").Value = "model_9_7_18"
$ws.Range("B20").Value = 0.9020052881573387
$ws.Range("C20").Value = 0.7203864744623816
$ws.Range("D20").Value = 0.9332865512042831
$ws.Range("E20").Value = 0.8681822622065581
$ws.Range("F20").Value = 0.1084512695670128
$ws.Range("G20").Value = 0.2674002647399902
$ws.Range("H20").Value = 0.1628671139478683
$ws.Range("I20").Value = 0.218208059668541
$ws.Range("A21").Value = "model_9_7_19"
$ws.Range("B21").Value = 0.9020572058008475
$ws.Range("C21").Value = 0.7206109280800526
$ws.Range("D21").Value = 0.9333046551495421
$ws.Range("E21").Value = 0.8682637874675794
$ws.Range("F21").Value = 0.1083938106894493
$ws.Range("G21").Value = 0.2671856284141541
$ws.Range("H21").Value = 0.1628229320049286
$ws.Range("I21").Value = 0.2180731147527695
$ws.Range("A22").Value = "model_9_7_23"
$ws.Range("B22").Value = 0.9021922248315151
$ws.Range("C22").Value = 0.7211245895744522
$ws.Range("D22").Value = 0.9332116268223001
$ws.Range("E22").Value = 0.8683555606526273
$ws.Range("F22").Value = 0.1082443818449974
$ws.Range("G22").Value = 0.2666943967342377
$ws.Range("H22").Value = 0.163050040602684
$ws.Range("I22").Value = 0.2179211974143982
$ws.Range("A23").Value = "model_9_7_24"
$ws.Range("B23").Value = 0.9022066904487339
$ws.Range("C23").Value = 0.7213903928782636
$ws.Range("D23").Value = 0.9330106289311837
$ws.Range("E23").Value = 0.8682980723926207
$ws.Range("F23").Value = 0.1082283779978752
$ws.Range("G23").Value = 0.266440212726593
$ws.Range("H23").Value = 0.1635407209396362
$ws.Range("I23").Value = 0.2180163562297821
$ws.Range("A24").Value = "model_9_7_20"
$ws.Range("B24").Value = 0.9022417025627009
$ws.Range("C24").Value = 0.7196165866199804
$ws.Range("D24").Value = 0.9342273945103956
$ws.Range("E24").Value = 0.8685997595656819
$ws.Range("F24").Value = 0.1081896275281906
$ws.Range("G24").Value = 0.2681365609169006
$ws.Range("H24").Value = 0.1605702489614487
$ws.Range("I24").Value = 0.2175169587135315
$ws.Range("A25").Value = "model_9_7_21"
$ws.Range("B25").Value = 0.9023302407047504
$ws.Range("C25").Value = 0.7212072764488655
$ws.Range("D25").Value = 0.9337649928233436
$ws.Range("E25").Value = 0.8687652727625801
$ws.Range("F25").Value = 0.1080916449427605
$ws.Range("G25").Value = 0.2666153013706207
$ws.Range("H25").Value = 0.1616991013288498
$ws.Range("I25").Value = 0.2172429859638214
$ws.Range("A26").Value = "model_9_7_22"
$ws.Range("B26").Value = 0.9023606893836169
$ws.Range("C26").Value = 0.721901239882354
$ws.Range("D26").Value = 0.9334766593612981
$ws.Range("E26").Value = 0.8687770591626337
$ws.Range("F26").Value = 0.1080579534173012
$ws.Range("G26").Value = 0.2659516930580139
$ws.Range("H26").Value = 0.162403017282486
$ws.Range("I26").Value = 0.2172234654426575
